$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellAddr, $value) {
    $rng = $ws.Range($cellAddr)
    $rng.NumberFormat = "@"
    $rng.Value = $value
    $rng.Style = "Normal"
}

Set-TextValue "D2" "29.938.68"
Set-TextValue "E2" "  +1.49%  "
Set-TextValue "D3" "1.939.90"
Set-TextValue "E3" "  +1.00%  "
Set-TextValue "D4" "1.009"
Set-TextValue "D5" "335.84"
Set-TextValue "E5" "  +2.91%  "
Set-TextValue "D6" "1.007"
Set-TextValue "E6" "  -0.07%  "
Set-TextValue "D7" "0.4840"
Set-TextValue "E7" "  +0.16%  "
Set-TextValue "D8" "0.4128"
Set-TextValue "E8" "  +1.02%  "
Set-TextValue "D9" "0.08191"
Set-TextValue "E9" "  -0.48%  "
Set-TextValue "D10" "1.016"
Set-TextValue "E10" "  -0.79%  "
Set-TextValue "D11" "23.83"
Set-TextValue "E11" "  +1.23%  "
Set-TextValue "D12" "1.942.97"
Set-TextValue "E12" "  +0.76%  "
Set-TextValue "D13" "6.100"
Set-TextValue "E13" "  +0.80%  "
Set-TextValue "D14" "7.312"
Set-TextValue "E14" "  +0.98%  "
Set-TextValue "D15" "91.38"
Set-TextValue "E15" "  -0.02%  "
Set-TextValue "D16" "0.06859"
Set-TextValue "E16" "  +0.72%  "
Set-TextValue "D17" "1.008"
Set-TextValue "E17" "  -0.07%  "
Set-TextValue "D18" "0.00001038"
Set-TextValue "E18" "  -0.36%  "
Set-TextValue "D19" "17.86"
Set-TextValue "E19" "  +0.15%  "
Set-TextValue "D20" "1.008"
Set-TextValue "E20" "  -0.03%  "
Set-TextValue "D21" "29.941.99"
Set-TextValue "E21" "  +1.42%  "
Set-TextValue "D22" "5.647"
Set-TextValue "E22" "  -0.05%  "
Set-TextValue "D23" "11.89"
Set-TextValue "E23" "  +0.92%  "
Set-TextValue "D24" "2.195"
Set-TextValue "E24" "  -0.05%  "
Set-TextValue "D25" "2.187.29"
Set-TextValue "E25" "  +1.00%  "
Set-TextValue "D26" "6.711"
Set-TextValue "E26" "  +0.37%  "
Set-TextValue "D27" "156.93"
Set-TextValue "E27" "  +0.03%  "
Set-TextValue "D28" "20.12"
Set-TextValue "E28" "  +0.11%  "
Set-TextValue "D29" "2.106"
Set-TextValue "E29" "  -0.68%  "
Set-TextValue "D30" "121.43"
Set-TextValue "E30" "  +0.70%  "
Set-TextValue "D31" "1.015"
Set-TextValue "E31" "  -1.23%  "
Set-TextValue "D32" "0.09638"
Set-TextValue "E32" "  +0.58%  "
Set-TextValue "D33" "5.607"
Set-TextValue "E33" "  +1.15%  "
Set-TextValue "D34" "1.423"
Set-TextValue "E34" "  +2.74%  "
Set-TextValue "D35" "3.549"
Set-TextValue "E35" "  -0.41%  "
Set-TextValue "D36" "0.06579"
Set-TextValue "E36" "  +7.12%  "
Set-TextValue "D37" "0.02291"
Set-TextValue "E37" "  +0.14%  "
Set-TextValue "D38" "1.216"
Set-TextValue "E38" "  +3.06%  "
Set-TextValue "D39" "0.5973"
Set-TextValue "E39" "  -0.39%  "
Set-TextValue "D40" "8.022"
Set-TextValue "E40" "  -0.41%  "
Set-TextValue "D41" "10.74"
Set-TextValue "E41" "  -0.89%  "
Set-TextValue "B42" "Algorand"
Set-TextValue "C42" "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
Set-TextValue "D42" "0.1855"
Set-TextValue "E42" "  -0.16%  "
Set-TextValue "B43" "RenderToken"
Set-TextValue "C43" "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
Set-TextValue "D43" "2.507"
Set-TextValue "E43" "  +4.11%  "
Set-TextValue "B44" "EnergySwap"
Set-TextValue "C44" "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
Set-TextValue "D44" "12.44"
Set-TextValue "E44" "  -0.06%  "
Set-TextValue "B45" "WEMIXToken"
Set-TextValue "C45" "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
Set-TextValue "D45" "1.238"
Set-TextValue "E45" "  -3.35%  "
Set-TextValue "D46" "0.07514"
Set-TextValue "E46" "  -1.20%  "
Set-TextValue "D47" "0.5576"
Set-TextValue "E47" "  -0.22%  "
Set-TextValue "D48" "1.988"
Set-TextValue "E48" "  +1.33%  "
Set-TextValue "D49" "117.44"
Set-TextValue "E49" "  -0.50%  "
Set-TextValue "D50" "72.91"
Set-TextValue "E50" "  -0.01%  "
Set-TextValue "D51" "2.423"
